$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '91.692.21'
Set-TextValue $ws.Range("E2") '  +1.25%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.155.42'
Set-TextValue $ws.Range("E3") '  +1.66%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.07%  '

# Row 5
Set-TextValue $ws.Range("D5") '239.83'
Set-TextValue $ws.Range("E5") '  -1.28%  '

# Row 6
Set-TextValue $ws.Range("D6") '619.41'
Set-TextValue $ws.Range("E6") '  -1.21%  '

# Row 7
Set-TextValue $ws.Range("D7") '1.11'
Set-TextValue $ws.Range("E7") '  -1.60%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.389'
Set-TextValue $ws.Range("E8") '  +5.01%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.999'
Set-TextValue $ws.Range("E9") '  -0.11%  '

# Row 10
Set-TextValue $ws.Range("D10") '3.153.98'
Set-TextValue $ws.Range("E10") '  +14.62%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.739'
Set-TextValue $ws.Range("E11") '  -0.77%  '

# Row 12
Set-TextValue $ws.Range("E12") '  +0.34%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.0000254'
Set-TextValue $ws.Range("E13") '  +1.85%  '

# Row 14
Set-TextValue $ws.Range("D14") '34.96'
Set-TextValue $ws.Range("E14") '  -1.65%  '

# Row 15
Set-TextValue $ws.Range("D15") '5.57'
Set-TextValue $ws.Range("E15") '  +1.06%  '

# Row 16
Set-TextValue $ws.Range("D16") '91.346.57'
Set-TextValue $ws.Range("E16") '  +1.02%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.732.64'
Set-TextValue $ws.Range("E17") '  +1.51%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.165.34'
Set-TextValue $ws.Range("E18") '  +0.21%  '

# Row 19
Set-TextValue $ws.Range("D19") '3.71'
Set-TextValue $ws.Range("E19") '  -3.72%  '

# Row 20
Set-TextValue $ws.Range("D20") '14.97'
Set-TextValue $ws.Range("E20") '  +4.71%  '

# Row 21
Set-TextValue $ws.Range("D21") '5.89'
Set-TextValue $ws.Range("E21") '  +1.78%  '

# Row 22
Set-TextValue $ws.Range("D22") '456.67'
Set-TextValue $ws.Range("E22") '  +2.50%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.0000202'
Set-TextValue $ws.Range("E23") '  -4.04%  '

# Row 24
Set-TextValue $ws.Range("D24") '9.17'
Set-TextValue $ws.Range("E24") '  +1.64%  '

# Row 25
Set-TextValue $ws.Range("B25") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D25") '1.62'
Set-TextValue $ws.Range("E25") '  +61.70%  '

# Row 26
Set-TextValue $ws.Range("B26") 'NEARProtocol'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D26") '5.67'
Set-TextValue $ws.Range("E26") '  -4.73%  '

# Row 27
Set-TextValue $ws.Range("D27") '88.68'
Set-TextValue $ws.Range("E27") '  -4.72%  '

# Row 28
Set-TextValue $ws.Range("D28") '11.75'
Set-TextValue $ws.Range("E28") '  -2.94%  '

# Row 29
Set-TextValue $ws.Range("B29") 'WrappedeETH'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range("D29") '3.319.72'
Set-TextValue $ws.Range("E29") '  +1.60%  '

# Row 30
Set-TextValue $ws.Range("B30") 'Hedera'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D30") '0.147'
Set-TextValue $ws.Range("E30") '  +31.44%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.999'
Set-TextValue $ws.Range("E31") '  -0.09%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.227'
Set-TextValue $ws.Range("E32") '  +3.94%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.167'
Set-TextValue $ws.Range("E33") '  -6.23%  '

# Row 34
Set-TextValue $ws.Range("D34") '9.34'
Set-TextValue $ws.Range("E34") '  +0.17%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.176'
Set-TextValue $ws.Range("E35") '  +12.51%  '

# Row 36
Set-TextValue $ws.Range("D36") '26.27'
Set-TextValue $ws.Range("E36") '  -1.30%  '

# Row 37
Set-TextValue $ws.Range("D37") '7.45'
Set-TextValue $ws.Range("E37") '  -1.87%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.94'
Set-TextValue $ws.Range("E38") '  +0.81%  '

# Row 39
Set-TextValue $ws.Range("D39") '490.31'
Set-TextValue $ws.Range("E39") '  -1.28%  '

# Row 40
Set-TextValue $ws.Range("B40") 'Fetch.AI'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D40") '1.32'
Set-TextValue $ws.Range("E40") '  +1.96%  '

# Row 41
Set-TextValue $ws.Range("B41") 'MantraDAO'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws.Range("D41") '3.87'
Set-TextValue $ws.Range("E41") '  -11.55%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.441'
Set-TextValue $ws.Range("E42") '  +5.67%  '

# Row 43
Set-TextValue $ws.Range("D43") '3.37'
Set-TextValue $ws.Range("E43") '  -7.00%  '

# Row 44
Set-TextValue $ws.Range("D44") '22.15'
Set-TextValue $ws.Range("E44") '  +0.19%  '

# Row 45
Set-TextValue $ws.Range("E45") '  +0.00%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.707'
Set-TextValue $ws.Range("E46") '  +2.46%  '

# Row 47
Set-TextValue $ws.Range("B47") 'Monero'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D47") '156.57'
Set-TextValue $ws.Range("E47") '  -1.98%  '

# Row 48
Set-TextValue $ws.Range("B48") 'Stacks'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D48") '1.92'
Set-TextValue $ws.Range("E48") '  +0.29%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.35'
Set-TextValue $ws.Range("E49") '  +0.05%  '

# Row 50
Set-TextValue $ws.Range("D50") '4.43'
Set-TextValue $ws.Range("E50") '  -2.87%  '

# Row 51
Set-TextValue $ws.Range("B51") 'OKB'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D51") '44.06'
Set-TextValue $ws.Range("E51") '  -2.06%  '
